# Add a new "Save" column (H) to the sheet:
#  - H1: header "Save", formatted like the other header cells (bold, centered, bordered)
#  - H2: value 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the adjacent "sum" header cell (G1) onto H1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header text and data value
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
